# Fruta / hortaliza, semanal
# Insert a new weekly record at row 265 (pushing the existing rows 265-333
# down to 266-334) and populate it with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 265; this shifts rows 265-333
# down to 266-334 and grows the sheet dimension to A1:R334.
$ws.Rows.Item(265).Insert()

# Fill in the data for the newly inserted row 265.
$ws.Cells.Item(265, 1).Value = 4
$ws.Cells.Item(265, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(265, 3).Value = "Los Lagos"
$ws.Cells.Item(265, 4).Value2 = 44932
$ws.Cells.Item(265, 5).Value = 10
$ws.Cells.Item(265, 6).Value = 100112032
$ws.Cells.Item(265, 7).Value = "Zapallo italiano"
$ws.Cells.Item(265, 8).Value = "Sin especificar"
$ws.Cells.Item(265, 9).Value = "Primera"
$ws.Cells.Item(265, 10).Value = 200
$ws.Cells.Item(265, 11).Value = 12000
$ws.Cells.Item(265, 12).Value = 12000
$ws.Cells.Item(265, 13).Value = 12000
$ws.Cells.Item(265, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(265, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(265, 16).Value = 240
$ws.Cells.Item(265, 17).Value = 50
$ws.Cells.Item(265, 18).Value = "Hortaliza"
